# Capacity Supply Curve.xlsx - calibration update
#  1) Update CSC parameters for dispatchable resources on the
#     "CSC-CSCSoCECBiaSY" sheet (Share of Cost Effective Capacity Built in a
#     Single Year): rows for dispatchable plant types go from 0.1 -> 0.5
#     across all year columns (B:AE). Intermittent/renewable resource rows
#     (hydro, onshore/offshore wind, solar PV/thermal, biomass, geothermal)
#     and municipal solid waste / crude oil are left untouched.
#  2) Leave the "About" sheet active/selected when the workbook is saved
#     (matches the author re-selecting the About tab before committing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

$rowRanges = @(
    @(2, 5),
    @(12, 14),
    @(16, 17),
    @(19, 25)
)

foreach ($rr in $rowRanges) {
    $startRow = $rr[0]
    $endRow = $rr[1]
    $rangeAddress = "B" + $startRow + ":AE" + $endRow
    $ws.Range($rangeAddress).Value = 0.5
}

# Re-select the "About" tab as the active sheet before saving.
$about = $wb.Worksheets.Item("About")
$about.Activate()
